$p = $ppt.ActivePresentation

# Hide slide at position 1 (sldId 428) and slide at position 4 (sldId 443)
$p.Slides.Item(1).SlideShowTransition.Hidden = $true
$p.Slides.Item(4).SlideShowTransition.Hidden = $true
